$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.675.21'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.730.83'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '563.38'
$ws.Range('E5').Value = '  -2.14%  '
$ws.Range('D6').Value = '158.63'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.593'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').Value = '0.108'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +3.25%  '
$ws.Range('D11').Value = '5.60'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = '3.216.54'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '26.69'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').Value = '63.519.77'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '0.0000148'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = '2.735.59'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '12.28'
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').Value = '352.37'
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').Value = '6.54'
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').Value = '0.520'
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').Value = '64.06'
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = '8.30'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '0.0₃0896'
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = '1.95'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '7.15'
$ws.Range('E30').Value = '  +3.12%  '
$ws.Range('E31').Value = '  +9.44%  '
$ws.Range('D32').Value = '163.07'
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('D33').Value = '19.97'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.84'
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').Value = '1.46'
$ws.Range('E36').Value = '  +2.14%  '
$ws.Range('D37').Value = '1.78'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = '0.974'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '342.39'
$ws.Range('E39').Value = '  +5.78%  '
$ws.Range('D40').Value = '6.22'
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').Value = '4.04'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').Value = '38.40'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').Value = '21.62'
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('D44').Value = '20.87'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').Value = '0.0576'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').Value = '134.03'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '0.621'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.0997'
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0248'
$ws.Range('E49').Value = '  -1.72%  '
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('E51').Value = '  +0.25%  '
